# Slide 7: table placeholder gains 3 new data rows and is re-centred
# ("putting the table always in the middle").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table

# Add three new rows at the end of the table (Tab-from-last-cell workflow).
$tbl.Rows.Add() | Out-Null
$tbl.Rows.Add() | Out-Null
$tbl.Rows.Add() | Out-Null

$newRowsData = @(
    @("12", "A", "`$13"),
    @("10", "A", "`$11"),
    @("23", "A", "`$25")
)

$startRow = $tbl.Rows.Count - $newRowsData.Count
for ($i = 0; $i -lt $newRowsData.Count; $i++) {
    $rowIndex = $startRow + $i + 1
    $rowData = $newRowsData[$i]
    for ($c = 1; $c -le 3; $c++) {
        $cell = $tbl.Cell($rowIndex, $c)
        $cell.Shape.TextFrame.TextRange.Text = $rowData[$c - 1]
    }
}

# Re-centre the (now taller) table within its placeholder: PowerPoint keeps
# growing tables vertically centred instead of only growing downward.
$sh.Top = 90.1511803023622
$sh.Height = 262.8000031
